$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.41945335488241
$ws.Range("C2").Value = 0.107605083706125
$ws.Range("B3").Value = -0.0399294119806398
$ws.Range("C3").Value = 0.0802228258285359

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.01849548600067
$ws.Range("C2").Value = 0.122883242611463
$ws.Range("B3").Value = -1.01997434819274
$ws.Range("C3").Value = 0.0818857483036436

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.94481970743779
$ws.Range("C2").Value = 0.105783665680401
$ws.Range("B3").Value = 1.78244840612281
$ws.Range("C3").Value = 0.145920306625554

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.17249190145231
$ws.Range("C2").Value = 0.111546443989113
$ws.Range("B3").Value = -0.0268888286602143
$ws.Range("C3").Value = 0.00923898580002034

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0115788540394022
$ws.Range("B2").Value = -0.00318512887970389
$ws.Range("A3").Value = -0.00318512887970389
$ws.Range("B3").Value = 0.00643570178391561

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0151002913147078
$ws.Range("B2").Value = -0.00742396882898494
$ws.Range("A3").Value = -0.00742396882898494
$ws.Range("B3").Value = 0.00670527577524767

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0111901839247827
$ws.Range("B2").Value = 0.00854226654679157
$ws.Range("A3").Value = 0.00854226654679157
$ws.Range("B3").Value = 0.0212927358856957

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0124426091666163
$ws.Range("B2").Value = -0.000441079631820998
$ws.Range("A3").Value = -0.000441079631820998
$ws.Range("B3").Value = 0.0000853588586129775
